$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "LowHealthMusic" row (row 20). This shifts all subsequent
# rows up by one, which is the bulk of the structural change in this
# revision (the LowHealthMusic asset/event was cut).
$ws.Rows(20).Delete()

# --- Status updates -----------------------------------------------------
# A number of sound/event pairs have been finished since the last pass,
# so flip their "Status of Sound" / "Status of Event" cells from
# Incomplete to Complete (and recolor them green to match the sheet's
# conditional-style convention: green fill = Complete, red fill = Incomplete).

function Set-Complete($addr) {
    $c = $ws.Range($addr)
    $c.Value = "Complete"
    $c.Interior.Color = 5296274
    $c.Font.Color = -16776961
}

# TutorialDialog (row 13)
Set-Complete("E13")
Set-Complete("F13")

# SpookyMusic (row 18)
Set-Complete("E18")
Set-Complete("F18")

# BossMusic (row 19)
Set-Complete("E19")
Set-Complete("F19")

# MenuMusic (row 20, formerly row 21 before the delete shifted it up)
Set-Complete("E20")
Set-Complete("F20")

# EnemyAmbiance (row 21, formerly row 22)
Set-Complete("E21")
Set-Complete("F21")

# PlayerAmbiance (row 22, formerly row 23)
Set-Complete("E22")
Set-Complete("F22")

# Exit Sign Buzz (row 24, formerly row 25)
Set-Complete("E24")
Set-Complete("F24")

# Staff Fire Crackling (row 25, formerly row 26) - status text updated but
# the original commit left the red Incomplete fill/font in place here.
$ws.Range("E25").Value = "Complete"
$ws.Range("F25").Value = "Complete"

# Staff Fire Crackling's categorization moved from Ambiance to Interface.
$ws.Range("C25").Value = "Interface"

# --- Notes updates -------------------------------------------------------
# TutorialDialog note rewritten now that status changed.
$ws.Range("H13").Value = "May need to re-record with stereo eventually. Didn't realize it was only recording mono till it was already in FMOD and I'm not sure how much I like it"

# The generic "not sure where I'm getting music..." placeholder note is no
# longer needed now that the music events are complete.
$ws.Range("H18").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("H20").Value = ""
